$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update Maria's balance (row 2, column C)
$ws.Range("C2").Value = 267446.15

# 2. Replace the Bluemetrix row (row 4) with Gustavo's data
# (leading apostrophe forces text so the leading zeros in the account
# number are preserved instead of being parsed away as a number)
$ws.Range("A4").Value = "'004444605"
$ws.Range("B4").Value = "Gustavo"
$ws.Range("C4").Value = 50700

# 3. Delete the Sergio row (account 004975924) entirely, shifting rows up
$ws.Rows(45).Delete()
